# After the team review, change some error code define.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the affected error-code cells in column A.
$ws.Range("A19").Value = "SKUM"
$ws.Range("A21").Value = "ERBN"
$ws.Range("A40").Value = "LWPS"
$ws.Range("A53").Value = "PTIN"

# Move the active selection from D4 to D3, as recorded in the sheet view.
$ws.Range("D3").Select()
